$wb = $excel.ActiveWorkbook

# --- Sheet "Feuil1" updates ---
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Range("C7").Value  = "coupon phoenix"
$ws1.Range("C8").Value  = "année"
$ws1.Range("C9").Value  = "12"
$ws1.Range("C10").Value = "Bouygues SA et BNP Paribas"
$ws1.Range("C11").Value = "wo action"
$ws1.Range("C13").Value = "degressif"
$ws1.Range("C14").Value = "1"
$ws1.Range("C15").Value = "81"
$ws1.Range("C16").Value = "80"
$ws1.Range("C17").Value = "80"
$ws1.Range("C24").Value = "29/07/2026"
$ws1.Range("C25").Value = "29/07/2027"
$ws1.Range("C28").Value = "02/08/2026"
$ws1.Range("C29").Value = "02/08/2027"

# --- Sheet "TRA" updates ---
$ws2 = $wb.Worksheets.Item("TRA")
$ws2.Range("A2").Value = "-22.13"
$ws2.Range("A3").Value = "1"
$ws2.Range("A4").Value = "1"
$ws2.Range("A5").Value = "1"
$ws2.Range("A6").Value = "1"

# --- Sheet "DATE" updates ---
$ws3 = $wb.Worksheets.Item("DATE")
$ws3.Range("A2").Value = "31/07/2023, 29/07/2024, 29/07/2025, 29/07/2026, 29/07/2027, 29/07/2027"
$ws3.Range("A3").Value = "Dates de paiement1"
$ws3.Range("A4").Value = "05/08/2022, 05/08/2022, 07/08/2023, 05/08/2024, 05/08/2025, 05/08/2026, 05/08/2027"
$ws3.Range("A5").Value = "Dates de remboursement"
$ws3.Range("A6").Value = "07/08/2023, 05/08/2024, 05/08/2025, 05/08/2026"
